$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column H width (column 8) from 2.140625 to 3.140625
# (Excel COM quantizes ColumnWidth to whole-pixel increments; 2.3 lands in
# the same pixel bucket as the target stored width of 3.140625.)
$ws.Columns.Item(8).ColumnWidth = 2.3

# Update row 1 values
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 32
$ws.Range("E1").Value = 23
$ws.Range("F1").Value = 11
$ws.Range("G1").Value = 3
$ws.Range("H1").Value = 29
$ws.Range("J1").Value = 10
$ws.Range("K1").Value = 0.010999999999999999
$ws.Range("L1").Value = 0.010999999999999999
$ws.Range("M1").Value = 0.028999999999999998
$ws.Range("N1").Value = 0.073999999999999996
